$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("B20").Value = "מורה יוכל להגביל את מספר המילים בתשובה"
